$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the 2021 data row (row 28), mirroring the existing table pattern ---

# Year
$ws.Range("A28").Value = 2021

# Raw counts
$ws.Range("E28").Value = 101
$ws.Range("F28").Value = 101
$ws.Range("G28").Value = 14

# Sightability-corrected adult counts (mirrors H22:I27 pattern: H=G, I=H)
$ws.Range("H28").Formula = "=G28"
$ws.Range("I28").Formula = "=H28"

# Adjusted counts
$ws.Range("J28").Formula = "=E28-H28"
$ws.Range("K28").Formula = "=F28-I28"

# Calves
$ws.Range("L28").Value = 37
$ws.Range("M28").Value = 37

# Sex ratio
$ws.Range("N28").Formula = "=M28/K28"

# Averages
$ws.Range("O28").Formula = "=AVERAGE(J28:K28)"
$ws.Range("P28").Formula = "=AVERAGE(H28:I28)"

# Sightability mean/SD
$ws.Range("Q28").Value = 1
$ws.Range("R28").Value = 0.09

# Estimates
$ws.Range("S28").Formula = "=AVERAGE(J28:K28)/Q28"
$ws.Range("T28").Formula = "=AVERAGE(J28:K28)*R28"
$ws.Range("U28").Formula = "=AVERAGE(H28:I28)/Q28"
$ws.Range("V28").Formula = "=AVERAGE(H28:I28)*R28"

# Citation
$ws.Range("W28").Value = "McNay Unpublished"

# --- View state: scroll/selection moved as part of editing the new row ---
$ws.Range("M29").Select()
